$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Duong Van Hien's cell: "-làm phần category" -> "-làm phần user"
#    (Do this FIRST, while "làm phần category" is still unique in the doc,
#    before we introduce a new "-làm phần category" paragraph elsewhere.)
# ---------------------------------------------------------------------------
$rngCategory = $d.Content
$foundCategory = $rngCategory.Find.Execute("-làm phần category", $true, $false, $false, $false, $false, `
                                            $true, 1, $false, "-làm phần user", 2)

# ---------------------------------------------------------------------------
# 2) Hoa Ngoc Anh's cell: remove the "-làm phần users" run, keep the (already
#    correctly formatted) empty paragraph that held it.
# ---------------------------------------------------------------------------
$rngUsers = $d.Content
$foundUsers = $rngUsers.Find.Execute("-làm phần users", $false, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)
if ($foundUsers) {
  $rngUsers.Text = ""
}

# ---------------------------------------------------------------------------
# 3) Le Quoc Anh's cell: remove the "-" + "làm trang detail" runs, keep the
#    (already correctly formatted) empty paragraph that held them.
# ---------------------------------------------------------------------------
$rngDetail = $d.Content
$foundDetail = $rngDetail.Find.Execute("-làm trang detail", $false, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)
if ($foundDetail) {
  $rngDetail.Text = ""
}

# ---------------------------------------------------------------------------
# 4) Tran Trung Kien's cell: after "...làm trang index, làm phần posts", add
#    three new paragraphs:
#       -làm phần category
#       -làm view
#       -tổng kết, sửa lỗi
#    Each one uses sz=32 / szCs=32 (16pt) like its siblings.
# ---------------------------------------------------------------------------
$newLines = @("-làm phần category", "-làm view", "-tổng kết, sửa lỗi")

$idx = 0
$postsIndex = -1
foreach ($par in $d.Paragraphs) {
  $idx++
  if ($par.Range.Text -like "*làm trang index, làm phần posts*") {
    $postsIndex = $idx
  }
}

if ($postsIndex -gt 0) {
  $idx = 0
  foreach ($par in $d.Paragraphs) {
    $idx++
    if ($idx -eq $postsIndex) {
      $par.Range.InsertParagraphAfter()
    }
  }

  $insertAt = $postsIndex + 1
  $lastLineIndex = $newLines.Count - 1
  for ($lineIdx = 0; $lineIdx -lt $newLines.Count; $lineIdx++) {
    $line = $newLines[$lineIdx]
    $idx = 0
    foreach ($par in $d.Paragraphs) {
      $idx++
      if ($idx -eq $insertAt) {
        $par.Range.Text = $line
        $par.Range.Font.Size = 16
        $par.Range.Font.SizeBi = 16
        if ($lineIdx -lt $lastLineIndex) {
          $par.Range.InsertParagraphAfter()
        }
      }
    }
    $insertAt++
  }
}
